# Applies the cell-value changes from the commit diff to the ARM/ALC/CUL/GSM/LTW/WVR sheets
# of the Pandaemonium_Profits workbook (currentAveragePrice / Leve profit columns).
$wb = $excel.ActiveWorkbook

function Set-Cells($ws, $values) {
    foreach ($addr in $values.Keys) {
        $v = $values[$addr]
        if ($null -eq $v) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $v
        }
    }
}

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-Cells $ws @{
    "H6" = 133
    "I6" = 133
    "K6" = 399
    "M6" = -287
    "H74" = 4914.273
    "I74" = 4697.5
    "J74" = 5038.143
    "K74" = 4697.5
    "L74" = 5038.143
    "M74" = -3761.5
    "N74" = -6910.143
    "H77" = 4914.273
    "I77" = 4697.5
    "J77" = 5038.143
    "K77" = 23487.5
    "L77" = 25190.715
    "M77" = -18807.5
    "N77" = -34550.715
    "H106" = 3194.4443
    "I106" = 2391.6667
    "J106" = 4800
    "K106" = 2391.6667
    "L106" = 4800
    "M106" = -1760.6667
    "N106" = -6062
    "H129" = 929.7458
    "J129" = 1001.1321
    "L129" = 3003.3963
    "N129" = -13003.3963
    "H138" = 5033.8604
    "J138" = 5741.1562
    "L138" = 17223.4686
    "N138" = -27503.4686
}

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-Cells $ws @{
    "H121" = $null
    "I121" = $null
    "J121" = $null
    "K121" = $null
    "L121" = $null
    "N121" = $null
    "H122" = $null
    "I122" = $null
    "J122" = $null
    "K122" = $null
    "L122" = $null
    "M122" = $null
    "N122" = $null
    "H123" = $null
    "I123" = $null
    "J123" = $null
    "K123" = $null
    "L123" = $null
    "N123" = $null
    "H124" = $null
    "I124" = $null
    "J124" = $null
    "K124" = $null
    "L124" = $null
    "H125" = $null
    "I125" = $null
    "J125" = $null
    "K125" = $null
    "L125" = $null
    "N125" = $null
    "H126" = $null
    "I126" = $null
    "J126" = $null
    "K126" = $null
    "L126" = $null
    "M126" = $null
    "H127" = $null
    "I127" = $null
    "J127" = $null
    "K127" = $null
    "L127" = $null
    "N127" = $null
    "H128" = $null
    "I128" = $null
    "J128" = $null
    "K128" = $null
    "L128" = $null
    "N128" = $null
    "H129" = $null
    "I129" = $null
    "J129" = $null
    "K129" = $null
    "L129" = $null
    "N129" = $null
    "H130" = $null
    "I130" = $null
    "J130" = $null
    "K130" = $null
    "L130" = $null
    "N130" = $null
    "H131" = $null
    "I131" = $null
    "J131" = $null
    "K131" = $null
    "L131" = $null
    "H132" = $null
    "I132" = $null
    "J132" = $null
    "K132" = $null
    "L132" = $null
    "M132" = $null
    "N132" = $null
    "H133" = $null
    "I133" = $null
    "J133" = $null
    "K133" = $null
    "L133" = $null
    "N133" = $null
    "H134" = $null
    "I134" = $null
    "J134" = $null
    "K134" = $null
    "L134" = $null
    "N134" = $null
    "H135" = $null
    "I135" = $null
    "J135" = $null
    "K135" = $null
    "L135" = $null
    "N135" = $null
    "H137" = $null
    "I137" = $null
    "J137" = $null
    "K137" = $null
    "L137" = $null
    "H138" = $null
    "I138" = $null
    "J138" = $null
    "K138" = $null
    "L138" = $null
    "N138" = $null
    "H139" = $null
    "I139" = $null
    "J139" = $null
    "K139" = $null
    "L139" = $null
    "N139" = $null
    "H140" = $null
    "I140" = $null
    "J140" = $null
    "K140" = $null
    "L140" = $null
    "N140" = $null
    "H141" = $null
    "I141" = $null
    "J141" = $null
    "K141" = $null
    "L141" = $null
    "N141" = $null
}

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-Cells $ws @{
    "H109" = 1299.1428
    "I109" = 1025.6666
    "J109" = 2940
    "K109" = 3076.9998
    "L109" = 8820
    "M109" = -2036.9998
    "N109" = -10900
    "H112" = 3500
    "I112" = 2000
    "J112" = 3714.2856
    "K112" = 6000
    "L112" = 11142.8568
    "M112" = -4892
    "N112" = -13358.8568
    "H122" = 1117.9512
    "I122" = 893.3333
    "J122" = 1135.6842
    "K122" = 8039.9997
    "L122" = 10221.1578
    "M122" = -5589.9997
    "N122" = -15121.1578
    "H125" = 3000
    "I125" = 3000
    "J125" = 0
    "K125" = 9000
    "L125" = 0
    "M125" = $null
    "N125" = -4080
}

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-Cells $ws @{
    "H134" = 40514.445
    "J134" = 40514.445
    "L134" = 121543.335
    "N134" = -126613.335
}

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-Cells $ws @{
    "H133" = 50878.855
    "J133" = 50878.855
    "L133" = 50878.855
    "N133" = -55938.855
}

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-Cells $ws @{
    "H119" = 0
    "I119" = 0
    "J119" = 0
    "K119" = 0
    "L119" = 0
    "H120" = 69700
    "I120" = 0
    "J120" = 69700
    "K120" = 0
    "L120" = 69700
    "N120" = -79376
    "H121" = 70000
    "I121" = 0
    "J121" = 70000
    "K121" = 0
    "L121" = 70000
    "N121" = -73494
    "H122" = 1427.25
    "I122" = 1002
    "J122" = 1852.5
    "K122" = 3006
    "L122" = 5557.5
    "M122" = -556
    "N122" = -10457.5
    "H123" = 53017.4
    "I123" = 0
    "J123" = 53017.4
    "K123" = 0
    "L123" = 53017.4
    "N123" = -62817.4
    "H124" = 0
    "I124" = 0
    "J124" = 0
    "K124" = 0
    "L124" = 0
    "H125" = 56000
    "I125" = 0
    "J125" = 56000
    "K125" = 0
    "L125" = 56000
    "N125" = -65840
    "H126" = 1661.8823
    "I126" = 1696.8
    "J126" = 1400
    "K126" = 5090.4
    "L126" = 4200
    "M126" = -2620.4
    "N126" = -9140
    "H127" = 62500
    "I127" = 0
    "J127" = 62500
    "K127" = 0
    "L127" = 62500
    "N127" = -72420
    "H128" = 61800
    "I128" = 0
    "J128" = 61800
    "K128" = 0
    "L128" = 61800
    "N128" = -71760
    "H129" = 62500
    "I129" = 0
    "J129" = 62500
    "K129" = 0
    "L129" = 62500
    "N129" = -72500
    "H130" = 0
    "I130" = 0
    "J130" = 0
    "K130" = 0
    "L130" = 0
    "H131" = 48933.332
    "I131" = 0
    "J131" = 48933.332
    "K131" = 0
    "L131" = 48933.332
    "N131" = -59013.332
    "H132" = 2737.4517
    "I132" = 2910.875
    "J132" = 2142.8572
    "K132" = 8732.625
    "L132" = 6428.571599999999
    "M132" = -6202.625
    "N132" = -11488.5716
    "H133" = 62664
    "I133" = 0
    "J133" = 62664
    "K133" = 0
    "L133" = 62664
    "N133" = -72784
    "H135" = 43139
    "I135" = 0
    "J135" = 43139
    "K135" = 0
    "L135" = 43139
    "N135" = -53279
    "H136" = 3481.4656
    "I136" = 2980.6924
    "J136" = 3888.3438
    "K136" = 8942.0772
    "L136" = 11665.0314
    "M136" = -6392.0772
    "N136" = -16765.0314
    "H137" = 62380
    "I137" = 0
    "J137" = 62380
    "K137" = 0
    "L137" = 62380
    "N137" = -72580
    "H138" = 44660
    "I138" = 0
    "J138" = 44660
    "K138" = 0
    "L138" = 44660
    "N138" = -54940
    "H139" = 53141.668
    "I139" = 0
    "J139" = 53141.668
    "K139" = 0
    "L139" = 53141.668
    "N139" = -63421.668
    "H140" = 61214.5
    "I140" = 0
    "J140" = 61214.5
    "K140" = 0
    "L140" = 61214.5
    "N140" = -71574.5
    "H141" = 64949.5
    "I141" = 0
    "J141" = 64949.5
    "K141" = 0
    "L141" = 64949.5
    "N141" = -75309.5
}

